$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1): gene name column (M) plus 12 experiment columns (A:L).
# Write the M1 "gene" label first so it lands at shared-string index 20,
# then the exp1..exp12 labels across A1:L1 so they follow in order.
$ws.Range("M1").Value = "gene"

$expLabels = @("exp1","exp2","exp3","exp4","exp5","exp6","exp7","exp8","exp9","exp10","exp11","exp12")
$headerCols = @("A","B","C","D","E","F","G","H","I","J","K","L")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $expLabels[$i]
}

# The new A1:L1 header cells pick up a (blank) alignment style, matching the
# extra cellXfs entry introduced for the header row.
$ws.Range("A1:L1").HorizontalAlignment = 1

# Selection moves to H23 (below the data) after the edit.
$ws.Range("H23").Select() | Out-Null
